# Update the cached "Last Modified" date field shown on every slide layout
# and on the slide master (the field itself, type="datetimeFigureOut", is
# recalculated by PowerPoint at display time, but the cached text stored in
# the OOXML must also be updated so tools reading the static text see the
# new date).
$p = $ppt.ActivePresentation

$oldDate = "7/2/2018"
$newDate = "9/10/2018"

# Slide master.
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shape = $master.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# Every slide layout belonging to the master.
$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shape = $layout.Shapes.Item($i)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Bump the version number shown on the cover slide.
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shape = $slide1.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "Version 1.1") {
            $tr.Text = "Version 1.2"
        }
    }
}
